# "Fixes issues in the course student" - replace the word list with a new
# set of values and drop the now-unused formatting (hyperlink-style cell on
# A2, empty styled cells in column B, and the extra rows that are no longer
# needed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused helper column (B) entirely - shrinks the used range
# back down to column A.
$ws.Columns("B:B").Delete()

# Drop the trailing rows (6-10) that are no longer part of the data set.
$ws.Rows("6:10").Delete()

# New word list.
$ws.Range("A1").Value = "ABCDFSDF"
$ws.Range("A2").Value = "VDAG`$%^"
$ws.Range("A3").Value = "HCSYFGDY"
$ws.Range("A4").Value = "VSHGDHDDH"
$ws.Range("A5").Value = "BVSHVDGH"

# A2 previously carried a (now unused) hyperlink-like style; clear it so the
# cell goes back to the default formatting.
$ws.Range("A2").ClearFormats()

# The "Hyperlink" cell style is no longer referenced by any cell - remove it
# so it doesn't linger in the style table.
$wb.Styles("Hyperlink").Delete()

# Selection / active cell ends up on the last populated row, like it would
# after typing the final value in the list.
[void]$ws.Range("A5").Select()
